# Generate Report for Handoff
# The "71c0c1a2-4e8f-4603-bb96-32ebfbc8c149.md" file is now ready for handoff
# (was previously "Handed back: in sync with en-US"), and its Latest Handoff
# Datetime is refreshed on the zh-cn / de-de status sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("D2").Value = "2016-02-23 09:29:51"
$zhcn.Range("D3").Value = "2016-02-23 09:29:51"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("D2").Value = "2016-02-23 09:30:02"
$dede.Range("D3").Value = "2016-02-23 09:30:02"
